# Show & edit student grade:
#  - Rename the "Point" column header to "Fullname" (B1)
#  - Move the active selection from B3 to B2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Fullname"
$ws.Range("B2").Select()
